# [PHOENIX - 5854] test data change and minor workflow fix for TL automation
#
# Updates the "approvalDetails" test-data sheet:
#  - Refreshes several TL (Trade License) approver test values so they use the
#    "<Name> [<Code>]" format (instead of the old "<Name>~<Code>" / "<Name> ~ <Code>"
#    formats), several with a trailing space as produced by the original edit.
#  - Adjusts which cell is currently selected/scrolled to in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order matters: new shared strings are appended in the order they are first
# used, so we set the D52:D56 cells before D51 to reproduce the exact shared
# string table ordering from the target workbook.
$ws.Range("D52").Value = "TLSanitoryInspectorOne [TL_SI_01] "
$ws.Range("D53").Value = "TLAMOHOne [TL_AMOH_10] "
$ws.Range("D54").Value = "TLMHO [TL_MHO_01] "
$ws.Range("D55").Value = "TLCMOH [TL_CMOH_01]"
$ws.Range("D56").Value = "TLSanitarySupervisorOne [TL_SS_01] "
$ws.Range("D51").Value = "TLCommissionerOne [ADM_COMM_1] "

# Minor workflow/view fix: move the active selection from D15 to D60.
$ws.Range("D60").Select()
